# Generate Report for Handoff
# Updates the "Latest Handoff Date"/"Latest Handoff Datetime" values for the
# 0280828d-6a5b-40e7-9911-09498559553f file row (row 5) across all three
# worksheets to reflect a freshly generated handoff report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D5").Value = "2016-42-13 08:42:11"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E5").Value = "2016-03-13 08:42:07"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E5").Value = "2016-03-13 08:42:11"
